$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# TestCase_F1 (row 2) result should now show SKIP, since only the D suite
# is being run (as noted in the commit message "running D suite only").
$ws.Range("D2").Value = "SKIP"
